# Update NATMI TPM-derived statistics for the Fgf1-Nrp1 ligand-receptor pair sheet.
# Each entry below is "<cell address>" = <new numeric value>, taken 1:1 from the
# updated TPM recomputation (see commit "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updatedValues = @{
    "G2" = 1.803092333333333
    "H2" = 5.409276999999999
    "I2" = 0.1744886524959502
    "J2" = 0.1744886524959502
    "M2" = 127.3992563333333
    "N2" = 382.197769
    "O2" = 0.4838549810199306
    "P2" = 0.4838549810199307
    "Q2" = 229.7126223670014
    "R2" = 2067.413601303013
    "S2" = 0.08442720364162128
    "T2" = 0.08442720364162129
    "G3" = 1.803092333333333
    "H3" = 5.409276999999999
    "I3" = 0.1744886524959502
    "J3" = 0.1744886524959502
    "M3" = 59.36586533333332
    "O3" = 0.2254681108101269
    "P3" = 0.2254681108101269
    "Q3" = 107.0421366442324
    "R3" = 963.3792297980918
    "S3" = 0.03934162683606663
    "T3" = 0.03934162683606663
    "G4" = 1.803092333333333
    "H4" = 5.409276999999999
    "I4" = 0.1744886524959502
    "J4" = 0.1744886524959502
    "M4" = 16.63275166666667
    "N4" = 49.898255
    "O4" = 0.06317022542837675
    "P4" = 0.06317022542837675
    "Q4" = 29.99038701240389
    "R4" = 269.9134831116349
    "S4" = 0.01102248751286287
    "T4" = 0.01102248751286287
    "G5" = 1.803092333333333
    "H5" = 5.409276999999999
    "I5" = 0.1744886524959502
    "J5" = 0.1744886524959502
    "M5" = 59.90262233333334
    "N5" = 179.707867
    "O5" = 0.2275066827415657
    "P5" = 0.2275066827415658
    "Q5" = 108.0099590757954
    "R5" = 972.0896316821591
    "S5" = 0.03969733450539947
    "T5" = 0.03969733450539947
    "I6" = 0.4384883998568034
    "J6" = 0.4384883998568034
    "M6" = 127.3992563333333
    "N6" = 382.197769
    "O6" = 0.4838549810199306
    "P6" = 0.4838549810199307
    "Q6" = 577.2657348646462
    "R6" = 5195.391613781816
    "S6" = 0.2121647963901734
    "T6" = 0.2121647963901734
    "I7" = 0.4384883998568034
    "J7" = 0.4384883998568034
    "M7" = 59.36586533333332
    "O7" = 0.2254681108101269
    "P7" = 0.2254681108101269
    "Q7" = 268.9959177458382
    "S7" = 0.09886515112786896
    "T7" = 0.09886515112786898
    "I8" = 0.4384883998568034
    "J8" = 0.4384883998568034
    "M8" = 16.63275166666667
    "N8" = 49.898255
    "O8" = 0.06317022542837675
    "P8" = 0.06317022542837675
    "Q8" = 75.36557033392445
    "R8" = 678.2901330053199
    "S8" = 0.02769941106668248
    "T8" = 0.02769941106668247
    "I9" = 0.4384883998568034
    "J9" = 0.4384883998568034
    "M9" = 59.90262233333334
    "N9" = 179.707867
    "O9" = 0.2275066827415657
    "P9" = 0.2275066827415658
    "Q9" = 271.4280467312543
    "R9" = 2442.852420581288
    "S9" = 0.0997590412720786
    "T9" = 0.09975904127207862
    "G10" = 3.895605666666667
    "H10" = 11.686817
    "I10" = 0.3769851220961256
    "J10" = 0.3769851220961256
    "M10" = 127.3992563333333
    "N10" = 382.197769
    "O10" = 0.4838549810199306
    "P10" = 0.4838549810199307
    "Q10" = 496.2972649012526
    "R10" = 4466.675384111273
    "S10" = 0.1824061290966171
    "T10" = 0.1824061290966171
    "G11" = 3.895605666666667
    "H11" = 11.686817
    "I11" = 0.3769851220961256
    "J11" = 0.3769851220961256
    "M11" = 59.36586533333332
    "O11" = 0.2254681108101269
    "P11" = 0.2254681108101269
    "Q11" = 231.2660013991035
    "R11" = 2081.394012591932
    "S11" = 0.08499812328253845
    "T11" = 0.08499812328253847
    "G12" = 3.895605666666667
    "H12" = 11.686817
    "I12" = 0.3769851220961256
    "J12" = 0.3769851220961256
    "M12" = 16.63275166666667
    "N12" = 49.898255
    "O12" = 0.06317022542837675
    "P12" = 0.06317022542837675
    "Q12" = 64.79464164492612
    "R12" = 583.1517748043351
    "S12" = 0.02381423514595639
    "T12" = 0.02381423514595639
    "G13" = 3.895605666666667
    "H13" = 11.686817
    "I13" = 0.3769851220961256
    "J13" = 0.3769851220961256
    "M13" = 59.90262233333334
    "N13" = 179.707867
    "O13" = 0.2275066827415657
    "P13" = 0.2275066827415658
    "Q13" = 233.3569950099266
    "R13" = 2100.21295508934
    "S13" = 0.08576663457101368
    "T13" = 0.08576663457101369
    "G14" = 0.1037266666666667
    "H14" = 0.31118
    "I14" = 0.01003782555112075
    "J14" = 0.01003782555112075
    "M14" = 127.3992563333333
    "N14" = 382.197769
    "O14" = 0.4838549810199306
    "P14" = 0.4838549810199307
    "Q14" = 13.21470019526889
    "R14" = 118.93230175742
    "S14" = 0.004856851891518906
    "T14" = 0.004856851891518906
    "G15" = 0.1037266666666667
    "H15" = 0.31118
    "I15" = 0.01003782555112075
    "J15" = 0.01003782555112075
    "M15" = 59.36586533333332
    "O15" = 0.2254681108101269
    "P15" = 0.2254681108101269
    "Q15" = 6.157823324808889
    "R15" = 55.42040992328
    "S15" = 0.002263209563652816
    "T15" = 0.002263209563652817
    "G16" = 0.1037266666666667
    "H16" = 0.31118
    "I16" = 0.01003782555112075
    "J16" = 0.01003782555112075
    "M16" = 16.63275166666667
    "N16" = 49.898255
    "O16" = 0.06317022542837675
    "P16" = 0.06317022542837675
    "Q16" = 1.725259887877778
    "R16" = 15.5273389909
    "S16" = 0.000634091702875018
    "T16" = 0.0006340917028750179
    "G17" = 0.1037266666666667
    "H17" = 0.31118
    "I17" = 0.01003782555112075
    "J17" = 0.01003782555112075
    "M17" = 59.90262233333334
    "N17" = 179.707867
    "O17" = 0.2275066827415657
    "P17" = 0.2275066827415658
    "Q17" = 6.21349933922889
    "R17" = 55.92149405306001
    "S17" = 0.002283672393074011
    "T17" = 0.002283672393074011
}

foreach ($cellAddress in $updatedValues.Keys) {
    $ws.Range($cellAddress).Value = $updatedValues[$cellAddress]
}
